$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 376, shifting existing rows 376-467 down to 377-468
$ws.Rows("376:376").Insert()

# Populate the newly inserted row 376 with the new data record
$ws.Range("A376").Value = 10
$ws.Range("B376").Value = "Vega Modelo de Temuco"
$ws.Range("C376").Value = "La Araucanía"
$ws.Range("D376").Value = 44722
$ws.Range("D376").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E376").Value = 9
$ws.Range("F376").Value = 100112032
$ws.Range("G376").Value = "Zapallo italiano"
$ws.Range("H376").Value = "Sin especificar"
$ws.Range("I376").Value = "Primera"
$ws.Range("J376").Value = 80
$ws.Range("K376").Value = 15000
$ws.Range("L376").Value = 15000
$ws.Range("M376").Value = 15000
$ws.Range("N376").Value = "$/caja 60 unidades"
$ws.Range("O376").Value = "Región de Arica y Parinacota"
$ws.Range("P376").Value = 250
$ws.Range("Q376").Value = 60
$ws.Range("R376").Value = "Hortaliza"
